# Updates the Price (column D) and Volume(1h) (column E) cells that changed
# in the latest cryptos-list scrape (GitHub Actions refresh).
#
# Every cell in these columns is stored as plain text in the workbook (e.g.
# "55.813.30", "  -1.70%  "), not a number. To stop Excel's COM layer from
# "smart"-converting numeric-looking literals (like "4.00" or "0.999") into
# real numbers, each value is written with a leading apostrophe (the same
# quote-prefix a person typing directly into Excel would use to force text).
# That quote-prefix nudges Excel into allocating a dedicated "quoted text"
# cell style, so immediately after the write we reset the cell back to the
# "Normal" style — this keeps the cell's formatting identical to the
# untouched cells (no stray numeric/text format is left behind) while the
# stored value itself stays exactly the literal text we wrote.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'55.813.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.70%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.346.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.87%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'504.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.02%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'129.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.52%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.00%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -2.57%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.356.03"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.77%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -0.34%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.01%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +2.09%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -1.29%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.761.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'55.759.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.63%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -0.37%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D18").Value = "'2.366.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.81%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -2.91%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'311.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.56%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.67%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.82%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.24%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'65.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.84%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.998"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.18%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.370"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.07%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -3.35%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -4.52%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'171.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.65%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.87%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -3.09%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.01%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -1.72%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.996"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.16%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -5.40%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.94%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -2.18%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.831"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.49%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -4.83%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -2.24%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -3.97%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -1.22%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'4.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.58%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'126.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.91%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -2.35%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0891"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.06%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'238.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.77%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0475"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.88%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -1.94%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'16.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.04%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.11%  "
$ws.Range("E51").Style = "Normal"
